$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove rows 3 and 4 (old fuzzy match rows no longer needed)
$ws.Range("A3:F4").Delete()

# Update row 2 with the new fuzzy-match results
$ws.Range("A2").Value = "PSD Cassidy Elementary School"
$ws.Range("B2").Value = "Kenderton Elementary School"
$ws.Range("C2").Value = "6523-43 LANSDOWNE AVE Lewis C. Cassidy Elementary School"
$ws.Range("D2").Value = "1500 W ONTARIO ST"
$ws.Range("E2").Value = 0.68
$ws.Range("F2").Value = "name"

# Adjust column widths (columns A & B narrower, C & D back to the sheet default)
$ws.Columns.Item(1).ColumnWidth = 20.166666666666668
$ws.Columns.Item(2).ColumnWidth = 13.5
$ws.Columns.Item(3).ColumnWidth = 8
$ws.Columns.Item(4).ColumnWidth = 8

# Reselect the active cell at A2
$ws.Range("A2").Select()
